$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 1289  # was 1284
$ws.Range("F5").Value = 276  # was 277
$ws.Range("F10").Value = 3542  # was 3536
$ws.Range("F13").Value = 72  # was 71
$ws.Range("F19").Value = 213  # was 212
$ws.Range("F22").Value = 66  # was 65
$ws.Range("F24").Value = 2742  # was 2728
$ws.Range("F25").Value = 5238  # was 5231
$ws.Range("F28").Value = 480  # was 479
$ws.Range("F29").Value = 3092  # was 3089
$ws.Range("F31").Value = 2269  # was 2267
$ws.Range("F35").Value = 133  # was 131
$ws.Range("F36").Value = 183  # was 182
$ws.Range("F38").Value = 36  # was 35
$ws.Range("F39").Value = 466  # was 464
$ws.Range("F40").Value = 810  # was 811
$ws.Range("F44").Value = 41  # was 40
$ws.Range("F45").Value = 494  # was 492

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 1289  # was 1284
$ws.Range("F5").Value = 276  # was 277
$ws.Range("F10").Value = 3542  # was 3536
$ws.Range("F13").Value = 72  # was 71
$ws.Range("F20").Value = 213  # was 212
$ws.Range("F23").Value = 66  # was 65
$ws.Range("F25").Value = 2742  # was 2728
$ws.Range("F26").Value = 5238  # was 5232
$ws.Range("F29").Value = 480  # was 479
$ws.Range("F30").Value = 3092  # was 3089
$ws.Range("F32").Value = 2269  # was 2267
$ws.Range("F36").Value = 133  # was 131
$ws.Range("F37").Value = 183  # was 182
$ws.Range("F39").Value = 36  # was 35
$ws.Range("F40").Value = 466  # was 464
$ws.Range("F41").Value = 810  # was 811
$ws.Range("F45").Value = 41  # was 40
$ws.Range("F46").Value = 494  # was 492

